# Fruta / hortaliza, semanal
# Updates weekly price/volume data for the Alcachofa (artichoke) subset,
# rotating dates, volumes, prices and origins across rows 2-16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44425
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("P2").Value = 467

# Row 3
$ws.Range("D3").Value = 44446
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 14000
$ws.Range("P3").Value = 467

# Row 4
$ws.Range("D4").Value = 44453
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 12000
$ws.Range("P4").Value = 400

# Row 5
$ws.Range("D5").Value = 44376
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 18000
$ws.Range("P5").Value = 600

# Row 6
$ws.Range("D6").Value = 44421
$ws.Range("J6").Value = 25
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15400
$ws.Range("P6").Value = 513

# Row 7
$ws.Range("D7").Value = 44435

# Row 8
$ws.Range("O8").Value = 'Provincia del Elquí'

# Row 9
$ws.Range("D9").Value = 44418
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 15000
$ws.Range("O9").Value = 'Provincia de Limarí'
$ws.Range("P9").Value = 500

# Row 10
$ws.Range("D10").Value = 44841
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("P10").Value = 400

# Row 11
$ws.Range("I11").Value = 'Segunda'
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("N11").Value = '$/caja 40 unidades'
$ws.Range("P11").Value = 250
$ws.Range("Q11").Value = 40

# Row 12
$ws.Range("D12").Value = 44474
$ws.Range("I12").Value = 'Primera'
$ws.Range("N12").Value = '$/caja 30 unidades'
$ws.Range("P12").Value = 333
$ws.Range("Q12").Value = 30

# Row 13
$ws.Range("D13").Value = 44432
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 14000
$ws.Range("O13").Value = 'Provincia del Elquí'
$ws.Range("P13").Value = 467

# Row 14
$ws.Range("D14").Value = 44449
$ws.Range("J14").Value = 45
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 12000
$ws.Range("P14").Value = 400

# Row 15
$ws.Range("D15").Value = 44467
$ws.Range("J15").Value = 35

# Row 16
$ws.Range("D16").Value = 44460
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 13000
$ws.Range("L16").Value = 13000
$ws.Range("M16").Value = 13000
$ws.Range("O16").Value = 'Provincia de Limarí'
$ws.Range("P16").Value = 433
